# "Basic load combinations generator"
#
# - Loadcases sheet: add two new Rail Track loadcase rows (12 & 13)
# - Envelopes sheet: rename the three single-loadcase envelope names to be
#   more specific, rename the repeated "Other" traffic envelope name to
#   "TLO Traffic Envelope", and add two new Rail Track envelope rows (15 & 16)
# - Widen/narrow a couple of the Envelopes columns to fit the new text
# - Make "Envelopes" the active sheet/tab (it was "Loadcases" before)

$wb = $excel.ActiveWorkbook

$wsLoadcases = $wb.Worksheets.Item("Loadcases")
$wsEnvelopes = $wb.Worksheets.Item("Envelopes")

# ---------------------------------------------------------------------------
# Envelopes sheet - more descriptive envelope names
# ---------------------------------------------------------------------------
$wsEnvelopes.Range("A2").Value = "Settlement Envelope"
$wsEnvelopes.Range("A3").Value = "Wind Envelope"

$wsEnvelopes.Range("A5").Value = "TLO Traffic Envelope"
$wsEnvelopes.Range("A6").Value = "TLO Traffic Envelope"
$wsEnvelopes.Range("A7").Value = "TLO Traffic Envelope"
$wsEnvelopes.Range("A8").Value = "TLO Traffic Envelope"
$wsEnvelopes.Range("A9").Value = "TLO Traffic Envelope"
$wsEnvelopes.Range("A10").Value = "TLO Traffic Envelope"
$wsEnvelopes.Range("A11").Value = "TLO Traffic Envelope"
$wsEnvelopes.Range("A12").Value = "TLO Traffic Envelope"
$wsEnvelopes.Range("A13").Value = "TLO Traffic Envelope"
$wsEnvelopes.Range("A14").Value = "TLO Traffic Envelope"

$wsEnvelopes.Range("A4").Value = "Thermal Envelope"

# ---------------------------------------------------------------------------
# Loadcases sheet - two new rows for rail track loadcases
# ---------------------------------------------------------------------------
$wsLoadcases.Range("B12").Value = 5
$wsLoadcases.Range("D12").Value = "Rail Load Characteristic"
$wsLoadcases.Range("B13").Value = 5
$wsLoadcases.Range("D13").Value = "Rail Load Characteristic"

# ---------------------------------------------------------------------------
# Envelopes sheet - two new rows for rail track envelopes
# ---------------------------------------------------------------------------
$wsEnvelopes.Range("A15").Value = "Track 1"

$wsLoadcases.Range("A12").Value = "Rail Track Up"
$wsLoadcases.Range("A13").Value = "Rail Track Down"

$wsEnvelopes.Range("B15").Value = "Rail Track Up"
$wsEnvelopes.Range("C15").Value = "Yes"

$wsEnvelopes.Range("A16").Value = "Track 2"
$wsEnvelopes.Range("B16").Value = "Rail Track Down"
$wsEnvelopes.Range("C16").Value = "Yes"

# ---------------------------------------------------------------------------
# Column widths - column A needs to be wider for the longer envelope names,
# column C can be narrower now the entries are short "Yes" values
# ---------------------------------------------------------------------------
$wsEnvelopes.Columns.Item(1).ColumnWidth = 24.0
$wsEnvelopes.Columns.Item(3).ColumnWidth = 15.083333333333334

# ---------------------------------------------------------------------------
# Selections + active sheet
# ---------------------------------------------------------------------------
$wsLoadcases.Range("A14").Select()
$wsEnvelopes.Range("D8").Select()
$wsEnvelopes.Activate()
